# Updated cryptos list on Wed Oct  4 00:27:51 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.393.43'
$ws.Range('E2').Value = '  -0.54%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.653.05'
$ws.Range('E3').Value = '  -0.72%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.78'
$ws.Range('E5').Value = '  -0.90%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.528'
$ws.Range('E6').Value = '  +3.31%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.51'
$ws.Range('E8').Value = '  +0.64%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.259'
$ws.Range('E9').Value = '  -0.17%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0611'
$ws.Range('E10').Value = '  -1.46%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0899'
$ws.Range('E11').Value = '  +2.56%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.889.30'
$ws.Range('E12').Value = '  -0.59%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.650.17'
$ws.Range('E13').Value = '  -0.80%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.06'
$ws.Range('E14').Value = '  -1.72%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.562'
$ws.Range('E15').Value = '  +2.37%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.97'
$ws.Range('E16').Value = '  -1.85%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.386.45'
$ws.Range('E17').Value = '  -0.58%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '230.16'
$ws.Range('E18').Value = '  -7.64%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0723'
$ws.Range('E19').Value = '  -1.28%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.39'
$ws.Range('E20').Value = '  -0.74%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  -0.06%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.32'
$ws.Range('E22').Value = '  -3.27%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.39'
$ws.Range('E23').Value = '  +0.66%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.03'
$ws.Range('E24').Value = '  +0.25%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.93'
$ws.Range('E25').Value = '  +0.51%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.03'
$ws.Range('E26').Value = '  -1.44%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.113'
$ws.Range('E27').Value = '  +1.58%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.11%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.61'
$ws.Range('E29').Value = '  -3.53%  '

# Row 30
$ws.Range('E30').Value = '  -3.81%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0491'
$ws.Range('E31').Value = '  -4.68%  '

# Row 32
$ws.Range('E32').Value = '  -1.12%  '

# Row 33
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.418.08'
$ws.Range('E33').Value = '  -2.55%  '

# Row 34
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.10'
$ws.Range('E34').Value = '  -0.93%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.56'
$ws.Range('E35').Value = '  +0.15%  '

# Row 36
$ws.Range('E36').Value = '  -0.02%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.896'
$ws.Range('E37').Value = '  -3.90%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.561'
$ws.Range('E38').Value = '  -2.57%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0166'
$ws.Range('E39').Value = '  -1.92%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.04'
$ws.Range('E40').Value = '  +0.09%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.01%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.51'
$ws.Range('E42').Value = '  +1.35%  '

# Row 43
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.22'
$ws.Range('E43').Value = '  +0.57%  '

# Row 44
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '64.57'
$ws.Range('E44').Value = '  -6.98%  '

# Row 45
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.798.31'
$ws.Range('E45').Value = '  -0.43%  '

# Row 46
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.788'
$ws.Range('E46').Value = '  -0.72%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.64'
$ws.Range('E47').Value = '  -3.77%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '87.52'
$ws.Range('E48').Value = '  -1.90%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').Value = '  -3.49%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0999'
$ws.Range('E50').Value = '  -1.45%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.70'
$ws.Range('E51').Value = '  -1.97%  '

